# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (header in G1) was regenerated from the underlying
# box-score data (strikeouts instead of the old "Strike#" derived metric).
# This applies the newly computed per-row K values to column G (rows 2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value, as produced by the regenerated
# save_data pipeline.
$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 0
    28 = 0
    29 = 1
    30 = 3
    31 = 1
    32 = 3
    33 = 0
    34 = 0
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 2
    40 = 1
    41 = 0
    42 = 0
    43 = 1
    44 = 0
    45 = 2
    46 = 2
    47 = 0
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 4
    56 = 1
    57 = 1
    58 = 2
    59 = 3
    60 = 0
    61 = 3
    62 = 3
    63 = 0
    64 = 1
    65 = 1
    66 = 1
    67 = 1
    68 = 1
    69 = 1
    70 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
